$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Several match rows were re-ordered upstream (the scraper re-sorted rows
#    that share the same kick-off timestamp). Columns A:E (index / pais /
#    torneio / temporada / data_partida) stay put; columns F:V (the actual
#    match + odds payload) swap between the two affected rows.
# ---------------------------------------------------------------------------
function Swap-RowPayload {
    param($row1, $row2, $firstCol, $lastCol)

    $vals1 = @{}
    $vals2 = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1[$c] = $ws.Cells.Item($row1, $c).Value2
        $vals2[$c] = $ws.Cells.Item($row2, $c).Value2
    }
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row1, $c).Value = $vals2[$c]
        $ws.Cells.Item($row2, $c).Value = $vals1[$c]
    }
}

# Column F = 6 ... Column V = 22
Swap-RowPayload 32 33 6 22
Swap-RowPayload 35 36 6 22
Swap-RowPayload 96 97 6 22
Swap-RowPayload 104 105 6 22

# ---------------------------------------------------------------------------
# 2) Three newly-scraped matches were appended at the bottom of the sheet
#    (rows 139-141). Clone the formatting of the last existing row (138) so
#    the index column keeps its bold/bordered style and the date column
#    keeps its date-time number format, then fill in the new data.
# ---------------------------------------------------------------------------
$ws.Range("A138:V138").Copy()
$ws.Range("A139:V141").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# "temporada" (column D) is a text column even though its content ("2023")
# looks numeric; force text storage for the three appended rows, then
# reapply D138's (style-less) formatting so the cell keeps the same
# (lack of) style as every other row while still storing a string value.
foreach ($r in 139..141) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = "2023"
}
$ws.Range("D138").Copy()
$ws.Range("D139:D141").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @{
        Row = 139; Index = 138; Date = 45225.97916666666
        Home = "San Lorenzo"; HomeGoals = 1; Away = "Platense"; AwayGoals = 1
        HomeOpenOdds = 1.93; HomeOpenDt = "19/10/2023 21:12"
        HomeCloseOdds = 2.11; HomeCloseDt = "26/10/2023 23:28"
        DrawOpenOdds = 2.99; DrawOpenDt = "19/10/2023 21:12"
        DrawCloseOdds = 2.77; DrawCloseDt = "26/10/2023 23:28"
        AwayOpenOdds = 5.25; AwayOpenDt = "19/10/2023 21:12"
        AwayCloseOdds = 5.04; AwayCloseDt = "26/10/2023 23:28"
        Url = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/san-lorenzo-platense/I5tVe67p/"
    },
    @{
        Row = 140; Index = 139; Date = 45226.08333333334
        Home = "Newells Old Boys"; HomeGoals = 0; Away = "Godoy Cruz"; AwayGoals = 2
        HomeOpenOdds = 2; HomeOpenDt = "20/10/2023 21:12"
        HomeCloseOdds = 2.15; HomeCloseDt = "27/10/2023 01:34"
        DrawOpenOdds = 3.16; DrawOpenDt = "20/10/2023 21:12"
        DrawCloseOdds = 3.04; DrawCloseDt = "27/10/2023 01:34"
        AwayOpenOdds = 4.07; AwayOpenDt = "20/10/2023 21:12"
        AwayCloseOdds = 4.18; AwayCloseDt = "27/10/2023 01:34"
        Url = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/newells-old-boys-godoy-cruz/dduZfQMj/"
    },
    @{
        Row = 141; Index = 140; Date = 45226.08333333334
        Home = "Lanus"; HomeGoals = 2; Away = "Tigre"; AwayGoals = 1
        HomeOpenOdds = 2.24; HomeOpenDt = "20/10/2023 21:12"
        HomeCloseOdds = 2.6; HomeCloseDt = "27/10/2023 01:59"
        DrawOpenOdds = 3.08; DrawOpenDt = "20/10/2023 21:12"
        DrawCloseOdds = 2.96; DrawCloseDt = "27/10/2023 01:28"
        AwayOpenOdds = 3.65; AwayOpenDt = "20/10/2023 21:12"
        AwayCloseOdds = 3.19; AwayCloseDt = "27/10/2023 01:56"
        Url = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/lanus-tigre/dYdn4Ua4/"
    }
)

foreach ($m in $newRows) {
    $r = $m.Row
    $ws.Cells.Item($r, 1).Value = $m.Index
    $ws.Cells.Item($r, 2).Value = "argentina"
    $ws.Cells.Item($r, 3).Value = "copa-de-la-liga-profesional"
    # column D ("2023") was already written further up so it's stored as
    # text instead of being auto-coerced to a number; don't touch it here.
    $ws.Cells.Item($r, 5).Value = $m.Date
    $ws.Cells.Item($r, 6).Value = $m.Home
    $ws.Cells.Item($r, 7).Value = $m.HomeGoals
    $ws.Cells.Item($r, 8).Value = $m.Away
    $ws.Cells.Item($r, 9).Value = $m.AwayGoals
    $ws.Cells.Item($r, 10).Value = $m.HomeOpenOdds
    $ws.Cells.Item($r, 11).Value = $m.HomeOpenDt
    $ws.Cells.Item($r, 12).Value = $m.HomeCloseOdds
    $ws.Cells.Item($r, 13).Value = $m.HomeCloseDt
    $ws.Cells.Item($r, 14).Value = $m.DrawOpenOdds
    $ws.Cells.Item($r, 15).Value = $m.DrawOpenDt
    $ws.Cells.Item($r, 16).Value = $m.DrawCloseOdds
    $ws.Cells.Item($r, 17).Value = $m.DrawCloseDt
    $ws.Cells.Item($r, 18).Value = $m.AwayOpenOdds
    $ws.Cells.Item($r, 19).Value = $m.AwayOpenDt
    $ws.Cells.Item($r, 20).Value = $m.AwayCloseOdds
    $ws.Cells.Item($r, 21).Value = $m.AwayCloseDt
    $ws.Cells.Item($r, 22).Value = $m.Url
}
